# Burndown diagram update: "new backlog" row (Sprint Reviews, row 20) is now
# included in the ACTUAL burndown calculation on row 23 (D23:L23), which
# previously only summed rows 4:19 and now sums rows 4:20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blatt1")

# Update the formulas in row 23 (ACTUAL burndown) so the daily burn sums
# include row 20 ("Sprint Reviews") in addition to rows 4-19.
$ws.Range("D23").Formula = "=C22-SUM(D4:D20)"
$ws.Range("E23").Formula = "=D23-SUM(E4:E20)"
$ws.Range("F23").Formula = "=E23-SUM(F4:F20)"
$ws.Range("G23").Formula = "=F23-SUM(G4:G20)"
$ws.Range("H23").Formula = "=G23-SUM(H4:H20)"
$ws.Range("I23").Formula = "=H23-SUM(I4:I20)"
$ws.Range("J23").Formula = "=I23-SUM(J4:J20)"
$ws.Range("K23").Formula = "=J23-SUM(K4:K20)"
$ws.Range("L23").Formula = "=K23-SUM(L4:L20)"

$excel.CalculateFullRebuild()

# Reflect the final selected cell in the saved sheet view.
$ws.Activate()
$ws.Range("M23").Select()

$wb.Save()
